# Update countries & provincias Spain
#
# The source data feed was refreshed (new case counts) and the sheet is kept
# sorted by "Casos totales" (column B) descending. That refresh changes the
# numbers for several countries and, because of the re-sort, swaps a few
# adjacent rows (El Salvador <-> Consejo Danes para los Refugiados,
# Estado de Palestina <-> Grecia, Uganda <-> Libia). The "Datos actualizados"
# timestamp in A1 is also bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "last updated" timestamp shown in A1.
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 11:13"

# Row => full A:H target values (country name, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes).
$data = @{
    37  = @("Ucrania", 47677, 914, 21115, 25335, 0, 15, 1227)
    39  = @("Singapur", 44664, 185, 39769, 4869, 0, 0, 26)
    40  = @("Oman", 43929, 0, 26169, 17557, 0, 10, 203)
    45  = @("Polonia", 35719, 314, 23127, 11080, 0, 5, 1512)
    60  = @("Austria", 18165, 115, 16607, 853, 0, 0, 705)
    74  = @("Malasia", 8658, 10, 8461, 76, 0, 0, 121)
    75  = @("Australia", 8362, 107, 7355, 903, 0, 0, 104)
    76  = @("El Salvador", 7507, 240, 4434, 2863, 0, 8, 210)
    77  = @("Consejo Danes para los Refugiados", 7311, 0, 2684, 4448, 0, 0, 179)
    98  = @("Estado de Palestina", 3599, 265, 463, 3125, 0, 0, 11)
    99  = @("Grecia", 3486, 0, 1374, 1920, 0, 0, 192)
    140 = @("Uganda", 927, 16, 849, 78, 0, 0, 0)
    141 = @("Libia", 918, 0, 230, 661, 0, 0, 27)
    187 = @("San Martin (Parte Holandesa)", 78, 1, 63, 0, 0, 0, 15)
    192 = @("Islas Turcas y Caicos", 45, 1, 11, 32, 0, 0, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
}
